$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells D1 / E1 -------------------------------------------------
$ws.Range("D1").Value = "porcentaje_utilidades"
$ws.Range("E1").Value = "porcentaje_contingencia"

# Start from the existing header formatting (bold font + centered/top alignment)
# by copying the format of an existing header cell, then trim the border down
# to left+right only (matching the new header style used for D1:E1).
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1:E1").Borders.Item(8).LineStyle = 0
$ws.Range("D1:E1").Borders.Item(9).LineStyle = 0
$excel.CutCopyMode = 0

# --- New data cell E2 ---------------------------------------------------------
$ws.Range("E2").Value = 13

# --- Column widths for C, D, E ------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 20.90625
$ws.Columns.Item(4).ColumnWidth = 20.08984375
$ws.Columns.Item(5).ColumnWidth = 23.36328125

# --- Selection matches the post-edit state in the source workbook ------------
$ws.Range("F2").Select()
